# Fruta / hortaliza, semanal
# Insert a new weekly record for "Macroferia Regional de Talca - Chirimoya"
# right before the existing row 46, shifting all the following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 46 (pushes old rows 46-66 down to 47-67)
$ws.Rows.Item(46).Insert()

# Populate the newly inserted row 46 with the new record's data
$ws.Cells.Item(46, 1).Value = 5
$ws.Cells.Item(46, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(46, 3).Value = "Maule"
$ws.Cells.Item(46, 4).Value = 44523
$ws.Cells.Item(46, 5).Value = 7
$ws.Cells.Item(46, 6).Value = "Fruta"
$ws.Cells.Item(46, 7).Value = 100107
$ws.Cells.Item(46, 8).Value = "Otros"
$ws.Cells.Item(46, 9).Value = 100107002
$ws.Cells.Item(46, 10).Value = "Chirimoya"
$ws.Cells.Item(46, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 150
$ws.Cells.Item(46, 14).Value = 23000
$ws.Cells.Item(46, 15).Value = 23000
$ws.Cells.Item(46, 16).Value = 23000
$ws.Cells.Item(46, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(46, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(46, 19).Value = 2300
$ws.Cells.Item(46, 20).Value = 10

# Make sure the date cell keeps the same date number format as the other
# rows in column D (numFmtId 165, style index 2 in the original workbook).
$ws.Cells.Item(46, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
